$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell, protecting against
# Excel auto-converting numeric-looking strings (e.g. "1.00", "7.70")
# into numbers. Values that are unambiguously non-numeric are written
# directly; ones that parse as a number get a quote-prefix so they stay
# text, then the style is reset back to Normal (no visible formatting).
function Set-TextCell($cell, $value, $forceText) {
    $c = $ws.Range($cell)
    if ($forceText) {
        $c.Value = "'" + $value
        $c.Style = "Normal"
    } else {
        $c.Value = $value
    }
}

Set-TextCell "D2" "64.034.62" $false
Set-TextCell "E2" "  -2.63%  " $false
Set-TextCell "D3" "2.632.60" $false
Set-TextCell "E3" "  -0.89%  " $false
Set-TextCell "E4" "  -0.06%  " $false
Set-TextCell "D5" "577.99" $true
Set-TextCell "E5" "  -3.27%  " $false
Set-TextCell "D6" "157.02" $true
Set-TextCell "E6" "  +0.03%  " $false
Set-TextCell "E7" "  -0.02%  " $false
Set-TextCell "E8" "  +0.49%  " $false
Set-TextCell "E9" "  -4.31%  " $false
Set-TextCell "E10" "  +0.50%  " $false
Set-TextCell "D11" "0.386" $true
Set-TextCell "E12" "  -0.32%  " $false
Set-TextCell "D13" "28.54" $true
Set-TextCell "D14" "3.108.64" $false
Set-TextCell "E14" "  -0.68%  " $false
Set-TextCell "E15" "  -5.63%  " $false
Set-TextCell "D16" "63.891.38" $false
Set-TextCell "E16" "  -2.57%  " $false
Set-TextCell "D17" "2.623.94" $false
Set-TextCell "E17" "  +0.36%  " $false
Set-TextCell "E18" "  -3.20%  " $false
Set-TextCell "D19" "7.70" $true
Set-TextCell "E19" "  +3.36%  " $false
Set-TextCell "D20" "4.63" $true
Set-TextCell "E20" "  -2.16%  " $false
Set-TextCell "D21" "345.77" $true
Set-TextCell "E21" "  -0.97%  " $false
Set-TextCell "E22" "  -0.15%  " $false
Set-TextCell "D23" "67.40" $true
Set-TextCell "E23" "  -2.50%  " $false
Set-TextCell "D24" "1.76" $true
Set-TextCell "E24" "  +1.87%  " $false
Set-TextCell "D25" "0.0000110" $true
Set-TextCell "E25" "  -1.55%  " $false
Set-TextCell "D26" "598.85" $true
Set-TextCell "E26" "  +10.04%  " $false
Set-TextCell "D27" "9.29" $true
Set-TextCell "E27" "  -3.22%  " $false
Set-TextCell "E28" "  +0.12%  " $false
Set-TextCell "E29" "  -0.72%  " $false
Set-TextCell "D30" "1.00" $true
Set-TextCell "E30" "  +0.01%  " $false
Set-TextCell "D31" "7.94" $true
Set-TextCell "E31" "  +0.54%  " $false
Set-TextCell "D32" "2.09" $true
Set-TextCell "E32" "  -2.51%  " $false
Set-TextCell "E33" "  -1.57%  " $false
Set-TextCell "D34" "6.64" $true
Set-TextCell "E34" "  +2.71%  " $false
Set-TextCell "D35" "5.35" $true
Set-TextCell "D36" "0.409" $true
Set-TextCell "E36" "  -1.96%  " $false
Set-TextCell "E37" "  -1.93%  " $false
Set-TextCell "E38" "  -0.01%  " $false
Set-TextCell "D39" "155.04" $true
Set-TextCell "E39" "  +0.11%  " $false
Set-TextCell "E40" "  -2.39%  " $false
Set-TextCell "E41" "  -0.01%  " $false
Set-TextCell "D42" "41.56" $true
Set-TextCell "E42" "  -2.15%  " $false
Set-TextCell "E43" "  +7.03%  " $false
Set-TextCell "D44" "157.70" $true
Set-TextCell "E44" "  -2.34%  " $false
Set-TextCell "E45" "  -2.28%  " $false
Set-TextCell "D46" "23.26" $true
Set-TextCell "E46" "  +3.12%  " $false
Set-TextCell "D47" "0.0600" $true
Set-TextCell "E47" "  -0.62%  " $false
Set-TextCell "B48" "Stellar" $false
Set-TextCell "C48" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm" $false
Set-TextCell "D48" "0.102" $true
Set-TextCell "E48" "  +2.44%  " $false
Set-TextCell "B49" "Mantle" $false
Set-TextCell "C49" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" $false
Set-TextCell "D49" "0.632" $true
Set-TextCell "E49" "  -0.86%  " $false
Set-TextCell "E50" "  -1.50%  " $false
Set-TextCell "D51" "19.09" $true
Set-TextCell "E51" "  -3.19%  " $false
